# Auto-generated edit script: refresh cryptos price/volume snapshot
# (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.126.34'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '2.373.38'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''303.90'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').Value = '''95.88'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').Value = '''0.503'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = '''1.00'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = '''0.483'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.62%  '
$ws.Range('D10').Value = '''34.40'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('D12').Value = '''0.0787'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').Value = '''18.56'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.37%  '
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('D15').Value = '2.735.79'
$ws.Range('E15').Value = '  +1.59%  '
$ws.Range('D16').Value = '2.382.11'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').Value = '''0.799'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('D18').Value = '43.071.90'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').Value = '''12.00'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('D20').Value = '''6.30'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.96%  '
$ws.Range('D21').Value = '0.0₃0888'
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('D22').Value = '''68.21'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.31%  '
$ws.Range('D23').Value = '''235.69'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.56%  '
$ws.Range('E24').Value = '  -2.60%  '
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = '''24.58'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('E28').Value = '  +15.22%  '
$ws.Range('D29').Value = '''9.39'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.66%  '
$ws.Range('D30').Value = '''32.18'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.69%  '
$ws.Range('D31').Value = '''0.999'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('D32').Value = '''5.01'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').Value = '''17.61'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('D34').Value = '''0.0718'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.09%  '
$ws.Range('D35').Value = '''0.106'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +4.76%  '
$ws.Range('E36').Value = '  +2.04%  '
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('E38').Value = '  +3.69%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').Value = '''124.28'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -10.28%  '
$ws.Range('B40').Value = 'WEMIXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').Value = '''2.26'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.65%  '
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('D42').Value = '''21.12'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -5.23%  '
$ws.Range('D43').Value = '1.936.42'
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('E45').Value = '  +3.93%  '
$ws.Range('D46').Value = '''9.32'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -7.33%  '
$ws.Range('D47').Value = '''2.73'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('D48').Value = '2.593.57'
$ws.Range('E48').Value = '  +1.27%  '
$ws.Range('E49').Value = '  +2.56%  '
$ws.Range('D50').Value = '''71.88'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('D51').Value = '''1.14'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.11%  '
